$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to Text format first so Excel does not
    # reinterpret a numeric-looking string (e.g. "0.999") as a number,
    # then clear the formatting back off so no extra style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "69.426.84"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "3.681.49"
$ws.Range("E3").Value = "  -2.33%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "682.49"
$ws.Range("E5").Value = "  -1.79%  "
Set-TextValue $ws.Range("D6") "159.57"
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("D7").Value = "3.680.70"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("E8").Value = "  -0.11%  "
Set-TextValue $ws.Range("D9") "0.493"
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("E10").Value = "  -8.37%  "
Set-TextValue $ws.Range("D11") "7.08"
$ws.Range("E11").Value = "  -5.35%  "
Set-TextValue $ws.Range("D12") "0.434"
$ws.Range("E12").Value = "  -8.05%  "
$ws.Range("E13").Value = "  -5.76%  "
$ws.Range("D14").Value = "4.298.73"
$ws.Range("E14").Value = "  -2.45%  "
Set-TextValue $ws.Range("D15") "32.34"
$ws.Range("E15").Value = "  -8.97%  "
$ws.Range("D16").Value = "3.676.70"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "69.379.55"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("E18").Value = "  -0.40%  "
Set-TextValue $ws.Range("D19") "15.78"
$ws.Range("E19").Value = "  -8.80%  "
$ws.Range("E20").Value = "  -9.11%  "
Set-TextValue $ws.Range("D21") "469.30"
$ws.Range("E21").Value = "  -8.21%  "
Set-TextValue $ws.Range("D22") "10.02"
$ws.Range("E22").Value = "  -2.77%  "
Set-TextValue $ws.Range("D23") "0.646"
$ws.Range("E23").Value = "  -8.57%  "
Set-TextValue $ws.Range("D24") "79.62"
$ws.Range("E24").Value = "  -4.10%  "
$ws.Range("D25").Value = "3.828.23"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -11.67%  "
$ws.Range("E28").Value = "  -12.04%  "
Set-TextValue $ws.Range("D29") "9.15"
$ws.Range("E29").Value = "  -9.28%  "
Set-TextValue $ws.Range("D30") "2.69"
$ws.Range("E30").Value = "  -7.75%  "
$ws.Range("E31").Value = "  -10.78%  "
Set-TextValue $ws.Range("D32") "2.01"
$ws.Range("E32").Value = "  -9.14%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "6.54"
$ws.Range("E34").Value = "  -9.47%  "
Set-TextValue $ws.Range("D35") "26.73"
$ws.Range("E35").Value = "  -7.11%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.654.25"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D37") "0.161"
$ws.Range("E37").Value = "  -4.12%  "
Set-TextValue $ws.Range("D38") "8.12"
$ws.Range("E38").Value = "  -11.33%  "
Set-TextValue $ws.Range("D39") "6.12"
$ws.Range("E39").Value = "  -5.32%  "
Set-TextValue $ws.Range("D40") "2.26"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("E42").Value = "  -9.15%  "
$ws.Range("E43").Value = "  -0.12%  "
Set-TextValue $ws.Range("D44") "0.941"
$ws.Range("E44").Value = "  -5.99%  "
Set-TextValue $ws.Range("D45") "165.10"
$ws.Range("E45").Value = "  +2.24%  "
Set-TextValue $ws.Range("D46") "47.48"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("E48").Value = "  -12.58%  "
Set-TextValue $ws.Range("D49") "1.30"
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("E50").Value = "  -7.01%  "
Set-TextValue $ws.Range("D51") "27.84"
$ws.Range("E51").Value = "  -4.79%  "
